$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new Rocket Launcher column
$ws.Range("C1").Value = "Rocket Launcher"

# Values for the new column (rows 2-13)
$values = @(0, 3, 1, 40, 20, 0.8, 1, 2, "-", 1, 1, 3)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Column width for C
$ws.Columns.Item(3).ColumnWidth = 14.6640625

# Alignment: right-align columns B and C (rows 1-13)
$ws.Range("B1:C13").HorizontalAlignment = -4152   # xlRight

# Border: medium gray border around top of data table (B1:C1 top/left/right)
$grayColor = 13421772   # RGB(204,204,204)

$ws.Range("B1:C1").Borders.Item(7).LineStyle = 1       # xlEdgeLeft
$ws.Range("B1:C1").Borders.Item(7).Weight = -4138      # xlMedium
$ws.Range("B1:C1").Borders.Item(7).Color = $grayColor

$ws.Range("B1:C1").Borders.Item(10).LineStyle = 1      # xlEdgeRight
$ws.Range("B1:C1").Borders.Item(10).Weight = -4138
$ws.Range("B1:C1").Borders.Item(10).Color = $grayColor

$ws.Range("B1:C1").Borders.Item(8).LineStyle = 1       # xlEdgeTop
$ws.Range("B1:C1").Borders.Item(8).Weight = -4138
$ws.Range("B1:C1").Borders.Item(8).Color = $grayColor

# Wrap text for C1 header
$ws.Range("C1").WrapText = $true

# Selection matches diff
$ws.Range("C13").Select()

# Page setup (paper size + orientation) as shown in diff
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait
